$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'254.56"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'3.42%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'28.09"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'-5.59%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'5.326"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'3.31%"
$ws.Range("E4").Style = "Normal"
$ws.Range("E5").Value = "'0.80%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'6.714"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'0.66%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.8666"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'1.72%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.9120"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'6.18%"
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.1425"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'3.69%"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.07172"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'1.03%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.03181"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'-0.68%"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.09221"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'-1.68%"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.001554"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'1.73%"
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.0006069"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'0.67%"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.005806"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'-2.03%"
$ws.Range("E15").Style = "Normal"
$ws.Range("E16").Value = "'0.06%"
$ws.Range("E16").Style = "Normal"
$ws.Range("E17").Value = "'0.05%"
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'2.242"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'1.39%"
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'0.3169"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'-0.88%"
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'0.03461"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'3.25%"
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'0.1315"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'1.46%"
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'3.552"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'1.90%"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.04168"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'0.71%"
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'0.1379"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'-0.14%"
$ws.Range("E24").Style = "Normal"
$ws.Range("E26").Value = "'-0.37%"
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "'-0.03%"
$ws.Range("E27").Style = "Normal"
$ws.Range("D40").Value = "'0.03840"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'2.23%"
$ws.Range("E40").Style = "Normal"
$ws.Range("E41").Value = "'2.78%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.002200"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'-9.49%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.002950"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'-16.43%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.01095"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'14.56%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.00005243"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'-0.92%"
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Value = "'0.01%"
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.08985"
$ws.Range("D47").Style = "Normal"
$ws.Range("E48").Value = "'-1.09%"
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'0.00002100"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'0.01%"
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.0002000"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'0.01%"
$ws.Range("E50").Style = "Normal"
